$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 164.83333
$ws.Range("I2").Value2 = 164.83333
$ws.Range("K2").Value2 = 164.83333
$ws.Range("M2").Value2 = -51.83332999999999
$ws.Range("H16").Value2 = 0
$ws.Range("J16").Value2 = 0
$ws.Range("L16").Value2 = 0
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value2 = 4022.0588
$ws.Range("J40").Value2 = 3648.4375
$ws.Range("L40").Value2 = 3648.4375
$ws.Range("N40").Value2 = -3998.4375
$ws.Range("H53").Value2 = 1055.2354
$ws.Range("J53").Value2 = 437.5
$ws.Range("L53").Value2 = 437.5
$ws.Range("N53").Value2 = -1711.5
$ws.Range("H62").Value2 = 5432.3335
$ws.Range("I62").Value2 = 5096.6
$ws.Range("K62").Value2 = 5096.6
$ws.Range("M62").Value2 = -4472.6
$ws.Range("H65").Value2 = 5432.3335
$ws.Range("I65").Value2 = 5096.6
$ws.Range("K65").Value2 = 25483
$ws.Range("M65").Value2 = -22363
$ws.Range("H70").Value2 = 5460
$ws.Range("J70").Value2 = 5460
$ws.Range("L70").Value2 = 16380
$ws.Range("N70").Value2 = -16920
$ws.Range("H73").Value2 = 5460
$ws.Range("J73").Value2 = 5460
$ws.Range("L73").Value2 = 16380
$ws.Range("N73").Value2 = -18252
$ws.Range("H100").Value2 = 1701.909
$ws.Range("I100").Value2 = 1780.25
$ws.Range("K100").Value2 = 1780.25
$ws.Range("M100").Value2 = -1239.25
$ws.Range("H132").Value2 = 6040.28
$ws.Range("I132").Value2 = 6229.4585
$ws.Range("K132").Value2 = 18688.3755
$ws.Range("M132").Value2 = -16158.3755
$ws.Range("H136").Value2 = 95902.664
$ws.Range("J136").Value2 = 100999.5
$ws.Range("L136").Value2 = 100999.5
$ws.Range("N136").Value2 = -111199.5
$ws.Range("H137").Value2 = 912134.9
$ws.Range("I137").Value2 = 1284203.9
$ws.Range("K137").Value2 = 3852611.7
$ws.Range("M137").Value2 = -3850061.7

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 1802388.1
$ws.Range("I32").Value2 = 843022.25
$ws.Range("J32").Value2 = 15873087
$ws.Range("K32").Value2 = 843022.25
$ws.Range("L32").Value2 = 15873087
$ws.Range("M32").Value2 = -842735.25
$ws.Range("N32").Value2 = -15873661
$ws.Range("H74").Value2 = 2134.0278
$ws.Range("I74").Value2 = 1692.238
$ws.Range("J74").Value2 = 2752.5334
$ws.Range("K74").Value2 = 1692.238
$ws.Range("L74").Value2 = 2752.5334
$ws.Range("M74").Value2 = -818.2380000000001
$ws.Range("N74").Value2 = -4500.5334
$ws.Range("H77").Value2 = 2134.0278
$ws.Range("I77").Value2 = 1692.238
$ws.Range("J77").Value2 = 2752.5334
$ws.Range("K77").Value2 = 8461.190000000001
$ws.Range("L77").Value2 = 13762.667
$ws.Range("M77").Value2 = -4093.190000000001
$ws.Range("N77").Value2 = -22498.667
$ws.Range("H97").Value2 = 790.7727
$ws.Range("I97").Value2 = 792.35
$ws.Range("J97").Value2 = 775
$ws.Range("K97").Value2 = 792.35
$ws.Range("L97").Value2 = 775
$ws.Range("M97").Value2 = -296.35
$ws.Range("N97").Value2 = -1767
$ws.Range("H102").Value2 = 3015.375
$ws.Range("I102").Value2 = 2386.2727
$ws.Range("K102").Value2 = 2386.2727
$ws.Range("M102").Value2 = -764.2727
$ws.Range("H132").Value2 = 4039.7778
$ws.Range("I132").Value2 = 2814.25
$ws.Range("J132").Value2 = 6490.8335
$ws.Range("K132").Value2 = 8442.75
$ws.Range("L132").Value2 = 19472.5005
$ws.Range("M132").Value2 = -5912.75
$ws.Range("N132").Value2 = -24532.5005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value2 = 0
$ws.Range("J76").Value2 = 0
$ws.Range("L76").Value2 = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value2 = 0
$ws.Range("J79").Value2 = 0
$ws.Range("L79").Value2 = 0
$ws.Range("N79").ClearContents()
$ws.Range("H99").Value2 = 93981
$ws.Range("I99").Value2 = 113476.78
$ws.Range("J99").Value2 = 6250
$ws.Range("K99").Value2 = 113476.78
$ws.Range("L99").Value2 = 6250
$ws.Range("M99").Value2 = -111978.78
$ws.Range("N99").Value2 = -9246
$ws.Range("H105").Value2 = 16251735
$ws.Range("I105").Value2 = 1430297
$ws.Range("K105").Value2 = 1430297
$ws.Range("M105").Value2 = -1428550

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3128317.2
$ws.Range("I31").Value2 = 1568.7391
$ws.Range("K31").Value2 = 1568.7391
$ws.Range("M31").Value2 = -1273.7391
$ws.Range("H34").Value2 = 3128317.2
$ws.Range("I34").Value2 = 1568.7391
$ws.Range("K34").Value2 = 1568.7391
$ws.Range("M34").Value2 = -1366.7391
$ws.Range("H59").Value2 = 76166.664
$ws.Range("J59").Value2 = 76166.664
$ws.Range("L59").Value2 = 76166.664
$ws.Range("N59").Value2 = -78456.664
$ws.Range("H132").Value2 = 4073.524
$ws.Range("I132").Value2 = 3424.2258
$ws.Range("J132").Value2 = 5903.364
$ws.Range("K132").Value2 = 10272.6774
$ws.Range("L132").Value2 = 17710.092
$ws.Range("M132").Value2 = -7742.6774
$ws.Range("N132").Value2 = -22770.092
$ws.Range("H134").Value2 = 3857.2632
$ws.Range("I134").Value2 = 4122.129
$ws.Range("K134").Value2 = 12366.387
$ws.Range("M134").Value2 = -9831.386999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value2 = 395
$ws.Range("J33").Value2 = 0
$ws.Range("L33").Value2 = 0
$ws.Range("N33").ClearContents()
$ws.Range("H34").Value2 = 570.4583
$ws.Range("J34").Value2 = 939.6
$ws.Range("L34").Value2 = 2818.8
$ws.Range("N34").Value2 = -2986.8
$ws.Range("H45").Value2 = 0
$ws.Range("J45").Value2 = 0
$ws.Range("L45").Value2 = 0
$ws.Range("N45").ClearContents()
$ws.Range("H80").Value2 = 3660.6667
$ws.Range("I80").Value2 = 2997
$ws.Range("J80").Value2 = 3992.5
$ws.Range("K80").Value2 = 8991
$ws.Range("L80").Value2 = 11977.5
$ws.Range("M80").Value2 = -8055
$ws.Range("N80").Value2 = -13849.5
$ws.Range("H83").Value2 = 3660.6667
$ws.Range("I83").Value2 = 2997
$ws.Range("J83").Value2 = 3992.5
$ws.Range("K83").Value2 = 26973
$ws.Range("L83").Value2 = 35932.5
$ws.Range("M83").Value2 = -22293
$ws.Range("N83").Value2 = -45292.5
$ws.Range("H97").Value2 = 1671833
$ws.Range("J97").Value2 = 7749.5
$ws.Range("L97").Value2 = 23248.5
$ws.Range("N97").Value2 = -24240.5
$ws.Range("H131").Value2 = 13893645
$ws.Range("I131").Value2 = 35723628
$ws.Range("J131").Value2 = 1838.4546
$ws.Range("K131").Value2 = 107170884
$ws.Range("L131").Value2 = 5515.3638
$ws.Range("M131").Value2 = -107165844
$ws.Range("N131").Value2 = -15595.3638
$ws.Range("H136").Value2 = 8576.182000000001
$ws.Range("I136").Value2 = 1223.6666
$ws.Range("J136").Value2 = 17399.2
$ws.Range("K136").Value2 = 3670.9998
$ws.Range("L136").Value2 = 52197.60000000001
$ws.Range("M136").Value2 = 1429.0002
$ws.Range("N136").Value2 = -62397.60000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 3963.3333
$ws.Range("J97").Value2 = 5000
$ws.Range("L97").Value2 = 5000
$ws.Range("N97").Value2 = -5992
$ws.Range("H105").Value2 = 54960.5
$ws.Range("J105").Value2 = 54960.5
$ws.Range("L105").Value2 = 54960.5
$ws.Range("N105").Value2 = -61948.5
$ws.Range("H126").Value2 = 7307.9443
$ws.Range("I126").Value2 = 2373.5
$ws.Range("J126").Value2 = 11255.5
$ws.Range("K126").Value2 = 7120.5
$ws.Range("L126").Value2 = 33766.5
$ws.Range("M126").Value2 = -4650.5
$ws.Range("N126").Value2 = -38706.5
$ws.Range("H141").Value2 = 26995
$ws.Range("J141").Value2 = 26995
$ws.Range("L141").Value2 = 26995
$ws.Range("N141").Value2 = -37355

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value2 = 9030
$ws.Range("I38").Value2 = 9030
$ws.Range("J38").Value2 = 0
$ws.Range("K38").Value2 = 9030
$ws.Range("L38").Value2 = 0
$ws.Range("M38").Value2 = -8620
$ws.Range("N38").ClearContents()
$ws.Range("H55").Value2 = 861.25
$ws.Range("I55").Value2 = 858
$ws.Range("K55").Value2 = 858
$ws.Range("M55").Value2 = -685
$ws.Range("H64").Value2 = 47997
$ws.Range("J64").Value2 = 47997
$ws.Range("L64").Value2 = 47997
$ws.Range("N64").Value2 = -48447
$ws.Range("H67").Value2 = 47997
$ws.Range("J67").Value2 = 47997
$ws.Range("L67").Value2 = 47997
$ws.Range("N67").Value2 = -49557
$ws.Range("H68").Value2 = 3203.4
$ws.Range("I68").Value2 = 3014
$ws.Range("K68").Value2 = 3014
$ws.Range("M68").Value2 = -2265
$ws.Range("H71").Value2 = 3203.4
$ws.Range("I71").Value2 = 3014
$ws.Range("K71").Value2 = 15070
$ws.Range("M71").Value2 = -11326
$ws.Range("H100").Value2 = 6358
$ws.Range("I100").Value2 = 5775
$ws.Range("J100").Value2 = 7524
$ws.Range("K100").Value2 = 5775
$ws.Range("L100").Value2 = 7524
$ws.Range("M100").Value2 = -5234
$ws.Range("N100").Value2 = -8606
$ws.Range("H132").Value2 = 4355.1875
$ws.Range("I132").Value2 = 4398.6924
$ws.Range("J132").Value2 = 4166.6665
$ws.Range("K132").Value2 = 13196.0772
$ws.Range("L132").Value2 = 12499.9995
$ws.Range("M132").Value2 = -10666.0772
$ws.Range("N132").Value2 = -17559.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value2 = 90910024
$ws.Range("I100").Value2 = 1183.8334
$ws.Range("J100").Value2 = 200000620
$ws.Range("K100").Value2 = 2367.6668
$ws.Range("L100").Value2 = 400001240
$ws.Range("M100").Value2 = -1826.6668
$ws.Range("N100").Value2 = -400002322
$ws.Range("H136").Value2 = 3719.2666
$ws.Range("I136").Value2 = 3870.7144
$ws.Range("K136").Value2 = 11612.1432
$ws.Range("M136").Value2 = -9062.143199999999
$ws.Range("H140").Value2 = 59885.332
$ws.Range("J140").Value2 = 59885.332
$ws.Range("L140").Value2 = 59885.332
$ws.Range("N140").Value2 = -70245.33199999999
$ws.Range("H141").Value2 = 69498.336
$ws.Range("J141").Value2 = 69498.336
$ws.Range("L141").Value2 = 69498.336
$ws.Range("N141").Value2 = -79858.336
